$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1,  "77C-226.75", "Bottom", "2025-11-02 18:25:25"),
    @(35, "77A-247.01", "Top",    "2025-11-02 18:25:41"),
    @(1,  "77C-226.75", "Bottom", "2025-11-18 22:40:17"),
    @(35, "77A-247.01", "Top",    "2025-11-18 22:40:44"),
    @(1,  "77C-226.75", "Bottom", "2025-11-18 22:51:28"),
    @(35, "77A-247.01", "Top",    "2025-11-18 22:51:48")
)

$row = 8
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}
